$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values in columns D and E stay as text
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '31.270.04'
$ws.Range("E2").Value = '  +2.98%  '
$ws.Range("D3").Value = '1.983.62'
$ws.Range("E3").Value = '  +6.01%  '
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '0.7950'
$ws.Range("E5").Value = '  +69.15%  '
$ws.Range("D6").Value = '252.92'
$ws.Range("E6").Value = '  +3.69%  '
$ws.Range("D7").Value = '0.9994'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '0.3376'
$ws.Range("E8").Value = '  +17.45%  '
$ws.Range("D9").Value = '25.70'
$ws.Range("E9").Value = '  +16.48%  '
$ws.Range("D10").Value = '0.06925'
$ws.Range("E10").Value = '  +7.38%  '
$ws.Range("D11").Value = '0.8331'
$ws.Range("E11").Value = '  +15.20%  '
$ws.Range("D12").Value = '0.08106'
$ws.Range("E12").Value = '  +4.47%  '
$ws.Range("D13").Value = '1.986.34'
$ws.Range("E13").Value = '  +6.21%  '
$ws.Range("D14").Value = '100.00'
$ws.Range("E14").Value = '  +4.24%  '
$ws.Range("D15").Value = '5.459'
$ws.Range("E15").Value = '  +6.41%  '
$ws.Range("D16").Value = '273.69'
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("D17").Value = '31.246.03'
$ws.Range("E17").Value = '  +2.95%  '
$ws.Range("D18").Value = '13.84'
$ws.Range("D19").Value = '0.000007935'
$ws.Range("E19").Value = '  +5.82%  '
$ws.Range("D22").Value = '0.9994'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '0.9995'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '6.921'
$ws.Range("D25").Value = '9.656'
$ws.Range("E25").Value = '  +6.71%  '
$ws.Range("D26").Value = '164.24'
$ws.Range("E26").Value = '  +0.68%  '
$ws.Range("D27").Value = '0.1468'
$ws.Range("E27").Value = '  +52.72%  '
$ws.Range("D28").Value = '19.76'
$ws.Range("E28").Value = '  +5.80%  '
$ws.Range("D29").Value = '2.168'
$ws.Range("E29").Value = '  +15.82%  '
$ws.Range("D30").Value = '1.564'
$ws.Range("E30").Value = '  +6.00%  '
$ws.Range("D31").Value = '1.352'
$ws.Range("E31").Value = '  +2.42%  '
$ws.Range("D32").Value = '4.553'
$ws.Range("E32").Value = '  +8.37%  '
$ws.Range("D33").Value = '4.327'
$ws.Range("E33").Value = '  +5.47%  '
$ws.Range("D34").Value = '0.05153'
$ws.Range("E34").Value = '  +7.30%  '
$ws.Range("D35").Value = '1.210'
$ws.Range("E35").Value = '  +8.21%  '
$ws.Range("D36").Value = '0.7544'
$ws.Range("E36").Value = '  +9.60%  '
$ws.Range("D37").Value = '2.765'
$ws.Range("E37").Value = '  +1.81%  '
$ws.Range("D38").Value = '0.9992'
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").Value = '0.02004'
$ws.Range("E39").Value = '  +7.00%  '
$ws.Range("D40").Value = '2.904'
$ws.Range("E40").Value = '  +3.38%  '
$ws.Range("D41").Value = '6.591'
$ws.Range("E41").Value = '  +6.36%  '
$ws.Range("D42").Value = '78.06'
$ws.Range("E42").Value = '  +5.22%  '
$ws.Range("D43").Value = '0.4645'
$ws.Range("E43").Value = '  +9.97%  '
$ws.Range("D44").Value = '2.052'
$ws.Range("E44").Value = '  +6.37%  '
$ws.Range("D45").Value = '0.8492'
$ws.Range("E45").Value = '  +2.45%  '
$ws.Range("D46").Value = '104.86'
$ws.Range("E46").Value = '  +4.04%  '
$ws.Range("D47").Value = '0.9992'
$ws.Range("D48").Value = '10.00'
$ws.Range("E48").Value = '  +4.87%  '
$ws.Range("D49").Value = '7.471'
$ws.Range("E49").Value = '  +7.68%  '
$ws.Range("D50").Value = '0.4283'
$ws.Range("E50").Value = '  +9.29%  '
$ws.Range("D51").Value = '36.47'
$ws.Range("E51").Value = '  +3.40%  '

# Row 20 and 21 content swap (Uniswap <-> WrappedliquidstakedEther2.0) with updated data
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '5.716'
$ws.Range("E20").Value = '  +9.44%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.247.79'
$ws.Range("E21").Value = '  +6.45%  '
